# Applies the "Added CDS All studies testcase" edit:
#  - The Samples query in B3 is narrowed to drop the Tumor / Analyte Type
#    columns (smp.sample_tumor_status / smp.sample_type).
#  - The Files query that used to live in B4 keeps its content unchanged.
#  - Selection/active cell moves from B4 to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$sampleQueryLines = @(
    'SELECT',
    '    DISTINCT (smp.sample_id) AS "Sample ID",',
    '    sp.participant_id AS "Participant ID", ',
    '    s.study_name AS "Study Name",',
    '    s.phs_accession AS Accession',
    'FROM ',
    '    df_participant sp',
    'JOIN ',
    '    df_study s ON sp."study.phs_accession" = s.phs_accession',
    'JOIN ',
    '    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id',
    'JOIN',
    '    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id',
    'JOIN',
    '    df_program p ON p.program_acronym = s."program.program_acronym"',
    'JOIN',
    '    df_file f1 ON f1."sample.sample_id" = smp.sample_id',
    'JOIN',
    '    df_genomic_info gi ON gi."file.file_id" = f1.file_id',
    'WHERE ',
    "   s.phs_accession = 'phs001524' AND gi.reference_genome_assembly = 'GRCh37'",
    'ORDER BY ',
    '    smp.sample_id ASC',
    'LIMIT 100;'
)
$newSampleQuery = [string]::Join("`r`n", $sampleQueryLines)

# Update the Samples query cell (B3) to the new, narrower query text.
$ws.Range("B3").Value = $newSampleQuery

# Move the selection to B3, matching the updated <selection> in the sheet view.
$ws.Range("B3").Select()
